# Apply data updates to Sheet1 reflecting refreshed survey calculations.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8: weight and number updated
$ws.Range("G8").Value = 4.075
$ws.Range("H8").Value = 50

# Rows 20-41 (except row 30 keeps its original Numb value):
# "Numb" (H) becomes -1 for all rows except row 30, and
# "RF" (I) is recalculated from 14.95096774193548 to 26.53967741935484
# for every row from 20 through 41.
$rf_new = 26.53967741935484

for ($r = 20; $r -le 41; $r++) {
    if ($r -ne 30) {
        $ws.Range("H$r").Value = -1
    }
    $ws.Range("I$r").Value = $rf_new
}
